$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2708.5
$ws.Range("I32").Value = 2633.9
$ws.Range("J32").Value = 2801.75
$ws.Range("K32").Value = 2633.9
$ws.Range("L32").Value = 2801.75
$ws.Range("M32").Value = -2307.9
$ws.Range("N32").Value = -3453.75

$ws.Range("H51").Value = 4266.6665
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 4200
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 4200
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -5168

$ws.Range("H81").Value = 34236
$ws.Range("J81").Value = 34236
$ws.Range("L81").Value = 34236
$ws.Range("N81").Value = -36232

$ws.Range("H84").Value = 34236
$ws.Range("J84").Value = 34236
$ws.Range("L84").Value = 102708
$ws.Range("N84").Value = -112692

$ws.Range("H98").Value = 2044.4062
$ws.Range("I98").Value = 1981.3226
$ws.Range("K98").Value = 1981.3226
$ws.Range("M98").Value = -483.3226

$ws.Range("H110").Value = 40561.6
$ws.Range("J110").Value = 40561.6
$ws.Range("L110").Value = 40561.6
$ws.Range("N110").Value = -48741.6

$ws.Range("H121").Value = 2903.2173
$ws.Range("J121").Value = 3008.1365
$ws.Range("L121").Value = 9024.4095
$ws.Range("N121").Value = -12518.4095

$ws.Range("H122").Value = 2044.4062
$ws.Range("I122").Value = 1981.3226
$ws.Range("K122").Value = 5943.9678
$ws.Range("M122").Value = -3493.9678

$ws.Range("H132").Value = 5311.8335
$ws.Range("I132").Value = 5880.1904
$ws.Range("J132").Value = 1333.3334
$ws.Range("K132").Value = 17640.5712
$ws.Range("L132").Value = 4000.0002
$ws.Range("M132").Value = -15110.5712
$ws.Range("N132").Value = -9060.0002

$ws.Range("H137").Value = 16667815
$ws.Range("I137").Value = 1133.0834
$ws.Range("J137").Value = 41667840
$ws.Range("K137").Value = 3399.2502
$ws.Range("L137").Value = 125003520
$ws.Range("M137").Value = -849.2501999999999
$ws.Range("N137").Value = -125008620

$ws.Range("H138").Value = 2319.5945
$ws.Range("I138").Value = 1747.738
$ws.Range("J138").Value = 3070.1562
$ws.Range("K138").Value = 5243.214
$ws.Range("L138").Value = 9210.4686
$ws.Range("M138").Value = -103.2139999999999
$ws.Range("N138").Value = -19490.4686

$ws.Range("H141").Value = 1404.4872
$ws.Range("I141").Value = 844.35486
$ws.Range("J141").Value = 3575
$ws.Range("K141").Value = 2533.06458
$ws.Range("L141").Value = 10725
$ws.Range("M141").Value = 2646.93542
$ws.Range("N141").Value = -21085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 866108.3
$ws.Range("I2").Value = 779.4761999999999
$ws.Range("J2").Value = 2263947.2
$ws.Range("K2").Value = 779.4761999999999
$ws.Range("L2").Value = 2263947.2
$ws.Range("M2").Value = -666.4761999999999
$ws.Range("N2").Value = -2264173.2

$ws.Range("H32").Value = 7353.395
$ws.Range("I32").Value = 5548.5073
$ws.Range("J32").Value = 15991.071
$ws.Range("K32").Value = 5548.5073
$ws.Range("L32").Value = 15991.071
$ws.Range("M32").Value = -5261.5073
$ws.Range("N32").Value = -16565.071

$ws.Range("H61").Value = 1372577
$ws.Range("I61").Value = 1588145.5
$ws.Range("J61").Value = 776.9091
$ws.Range("K61").Value = 1588145.5
$ws.Range("L61").Value = 776.9091
$ws.Range("M61").Value = -1587933.5
$ws.Range("N61").Value = -1200.9091

$ws.Range("H110").Value = 1287.3529
$ws.Range("I110").Value = 1098.7142
$ws.Range("J110").Value = 2167.6667
$ws.Range("K110").Value = 1098.7142
$ws.Range("L110").Value = 2167.6667
$ws.Range("M110").Value = 946.2858000000001
$ws.Range("N110").Value = -6257.6667

$ws.Range("H116").Value = 866108.3
$ws.Range("I116").Value = 779.4761999999999
$ws.Range("J116").Value = 2263947.2
$ws.Range("K116").Value = 779.4761999999999
$ws.Range("L116").Value = 2263947.2
$ws.Range("M116").Value = 1514.5238
$ws.Range("N116").Value = -2268535.2

$ws.Range("H122").Value = 1348
$ws.Range("I122").Value = 1303.6522
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 3910.9566
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -1460.9566
$ws.Range("N122").Value = -9250

$ws.Range("H132").Value = 5879362.5
$ws.Range("I132").Value = 7275949
$ws.Range("J132").Value = 78155.234
$ws.Range("K132").Value = 21827847
$ws.Range("L132").Value = 234465.702
$ws.Range("M132").Value = -21825317
$ws.Range("N132").Value = -239525.702

$ws.Range("H136").Value = 1372577
$ws.Range("I136").Value = 1588145.5
$ws.Range("J136").Value = 776.9091
$ws.Range("K136").Value = 4764436.5
$ws.Range("L136").Value = 2330.7273
$ws.Range("M136").Value = -4761886.5
$ws.Range("N136").Value = -7430.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 866108.3
$ws.Range("I3").Value = 779.4761999999999
$ws.Range("J3").Value = 2263947.2
$ws.Range("K3").Value = 779.4761999999999
$ws.Range("L3").Value = 2263947.2
$ws.Range("M3").Value = -665.4761999999999
$ws.Range("N3").Value = -2264175.2

$ws.Range("H94").Value = 933.8889
$ws.Range("I94").Value = 854.61536
$ws.Range("J94").Value = 1140
$ws.Range("K94").Value = 854.61536
$ws.Range("L94").Value = 1140
$ws.Range("M94").Value = -403.61536
$ws.Range("N94").Value = -2042

$ws.Range("H134").Value = 4372560
$ws.Range("I134").Value = 4571233
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 13713699
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -13711164
$ws.Range("N134").Value = -10320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5128411
$ws.Range("I31").Value = 1299.3334
$ws.Range("K31").Value = 1299.3334
$ws.Range("M31").Value = -1004.3334

$ws.Range("H34").Value = 5128411
$ws.Range("I34").Value = 1299.3334
$ws.Range("K34").Value = 1299.3334
$ws.Range("M34").Value = -1097.3334

$ws.Range("H58").Value = 2061.4922
$ws.Range("I58").Value = 1039.3784
$ws.Range("J58").Value = 3412.1428
$ws.Range("K58").Value = 1039.3784
$ws.Range("L58").Value = 3412.1428
$ws.Range("M58").Value = -836.3784000000001
$ws.Range("N58").Value = -3818.1428

$ws.Range("H132").Value = 2533.9312
$ws.Range("I132").Value = 2466.84
$ws.Range("J132").Value = 2953.25
$ws.Range("K132").Value = 7400.52
$ws.Range("L132").Value = 8859.75
$ws.Range("M132").Value = -4870.52
$ws.Range("N132").Value = -13919.75

$ws.Range("H134").Value = 17243238
$ws.Range("I134").Value = 2081.25
$ws.Range("J134").Value = 100000790
$ws.Range("K134").Value = 6243.75
$ws.Range("L134").Value = 300002370
$ws.Range("M134").Value = -3708.75
$ws.Range("N134").Value = -300007440

$ws.Range("H136").Value = 2061.4922
$ws.Range("I136").Value = 1039.3784
$ws.Range("J136").Value = 3412.1428
$ws.Range("K136").Value = 3118.1352
$ws.Range("L136").Value = 10236.4284
$ws.Range("M136").Value = -568.1352000000002
$ws.Range("N136").Value = -15336.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 33386.4
$ws.Range("J12").Value = 52682.58
$ws.Range("L12").Value = 158047.74
$ws.Range("N12").Value = -158393.74

$ws.Range("H111").Value = 1467.5
$ws.Range("I111").Value = 956.6667
$ws.Range("K111").Value = 2870.0001
$ws.Range("M111").Value = 196.9998999999998

$ws.Range("H113").Value = 448.48935
$ws.Range("I113").Value = 437.57895
$ws.Range("J113").Value = 465.2973
$ws.Range("K113").Value = 1312.73685
$ws.Range("L113").Value = 1395.8919
$ws.Range("M113").Value = 857.26315
$ws.Range("N113").Value = -5735.891900000001

$ws.Range("H122").Value = 8711849
$ws.Range("I122").Value = 15873577
$ws.Range("J122").Value = 1192035.5
$ws.Range("K122").Value = 142862193
$ws.Range("L122").Value = 10728319.5
$ws.Range("M122").Value = -142859743
$ws.Range("N122").Value = -10733219.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1589.2927
$ws.Range("I102").Value = 1491.7
$ws.Range("J102").Value = 1855.4546
$ws.Range("K102").Value = 1491.7
$ws.Range("L102").Value = 1855.4546
$ws.Range("M102").Value = 130.3
$ws.Range("N102").Value = -5099.4546

$ws.Range("H107").Value = 744
$ws.Range("I107").Value = 498.57144
$ws.Range("J107").Value = 958.75
$ws.Range("K107").Value = 498.57144
$ws.Range("L107").Value = 958.75
$ws.Range("M107").Value = 1421.42856
$ws.Range("N107").Value = -4798.75

$ws.Range("H132").Value = 23811410
$ws.Range("I132").Value = 35716190
$ws.Range("J132").Value = 1844.3572
$ws.Range("K132").Value = 107148570
$ws.Range("L132").Value = 5533.071599999999
$ws.Range("M132").Value = -107146040
$ws.Range("N132").Value = -10593.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5943.68
$ws.Range("I132").Value = 6955.2856
$ws.Range("J132").Value = 632.75
$ws.Range("K132").Value = 20865.8568
$ws.Range("L132").Value = 1898.25
$ws.Range("M132").Value = -18335.8568
$ws.Range("N132").Value = -6958.25

$ws.Range("H136").Value = 1254.5098
$ws.Range("I136").Value = 681.9737
$ws.Range("J136").Value = 2928.077
$ws.Range("K136").Value = 2045.9211
$ws.Range("L136").Value = 8784.231
$ws.Range("M136").Value = 504.0789
$ws.Range("N136").Value = -13884.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6130339
$ws.Range("I132").Value = 7005992.5
$ws.Range("J132").Value = 764.125
$ws.Range("K132").Value = 21017977.5
$ws.Range("L132").Value = 2292.375
$ws.Range("M132").Value = -21015447.5
$ws.Range("N132").Value = -7352.375

$ws.Range("H136").Value = 2980908.2
$ws.Range("I136").Value = 5625.3335
$ws.Range("J136").Value = 15873801
$ws.Range("K136").Value = 16876.0005
$ws.Range("L136").Value = 47621403
$ws.Range("M136").Value = -14326.0005
$ws.Range("N136").Value = -47626503
